# Merge the split "<id>...</id>" runs back into a single run for each of
# the two newly-downloaded <id> blocks (p100v_1 and p100v_2).
#
# In the original document each block is built from three separate runs:
#   <id>   (Courier New, color 7f6000, sz 18)
#   p100v_N  (separate run/formatting)
#   </id>  (Courier New, color 7f6000, sz 18)
#
# The edit collapses each triplet into one run whose text is the full
# "<id>p100v_N</id>" string, keeping the formatting of the first ("<id>")
# run. We achieve this with a literal Find & Replace across the whole
# document: replacing the already-contiguous visible text
# "<id>p100v_N</id>" with itself forces Word to re-write the matched
# range as a single run using the formatting of the run at the start of
# the match.

$d = $word.ActiveDocument

$ids = @("p100v_1", "p100v_2")

foreach ($id in $ids) {
    $text = "<id>" + $id + "</id>"

    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($text, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $text, 2)

    if (-not $found) {
        Write-Host "WARNING: pattern not found: $text"
    } else {
        Write-Host "Merged run for: $text"
    }
}
